$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 - new diary entry (5 Jan 2020)
$ws.Cells.Item(29, 1).Value = "5 Jan 2020 (W)"
$ws.Cells.Item(29, 2).Value = "1710-2010"
$ws.Cells.Item(29, 3).Value = "Harry, Deon, Thuc"
$ws.Cells.Item(29, 4).Value = "Work on fourth lecture’s homework while simutaneously paying attention to the 261 lecture"
$ws.Cells.Item(29, 5).Value = "We finished documenting the first feature"
$ws.Cells.Item(29, 6).Value = "Decisions, decisions! It’s hard to make a choice of features when there are so many! Also in Runeline a large chunk of the code is plugins. So the hunt for essential features became a lot harder."
$ws.Cells.Item(29, 7).Value = "Because the classes are in sequence, I have to work on the assignments in sequence. So the work for this class will always be at the tail end of the weekly schedule."

# Row 30 - continuation (5 Jan 2020)
$ws.Cells.Item(30, 1).Value = "5 Jan 2020 (W)"
$ws.Cells.Item(30, 2).Value = "2022-0000"
$ws.Cells.Item(30, 3).Formula = "=C29"
$ws.Cells.Item(30, 4).Formula = "=D29"
$ws.Cells.Item(30, 5).Value = "We finished documenting the second feature"

# Row 31 - new diary entry (6 Jan 2020)
$ws.Cells.Item(31, 1).Value = "6 Jan 2020 (Th)"
$ws.Cells.Item(31, 2).Value = "0000-0030"
$ws.Cells.Item(31, 3).Value = "Harry, Deon, Thuc"
$ws.Cells.Item(31, 4).Formula = "=D30"
$ws.Cells.Item(31, 5).Formula = "=E30"

$ws.Application.ActiveWindow.ScrollRow = 24
$ws.Range("B31").Select()
